$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: swap A5 and C5 values (서비스 ID <-> 인터페이스 ID)
$a5 = $ws.Range("A5").Value2
$c5 = $ws.Range("C5").Value2
$ws.Range("A5").Value = $c5
$ws.Range("C5").Value = $a5

# Row 7: reorder E7:J7
#   before: E7=어댑터 ID F7=인터페이스 ID G7=인스턴스 ID H7=메시지 ID I7=커넥터 ID J7=서비스 ID
#   after:  E7=메시지 ID F7=인터페이스 ID G7=서비스 ID  H7=인스턴스 ID I7=어댑터 ID  J7=커넥터 ID
$e7 = $ws.Range("E7").Value2
$g7 = $ws.Range("G7").Value2
$h7 = $ws.Range("H7").Value2
$i7 = $ws.Range("I7").Value2
$j7 = $ws.Range("J7").Value2

$ws.Range("E7").Value = $h7
$ws.Range("G7").Value = $j7
$ws.Range("H7").Value = $g7
$ws.Range("I7").Value = $e7
$ws.Range("J7").Value = $i7

# Clear the sheet view selection (reset to default / A1)
$ws.Range("A1").Select()
